# Appends one new day (2020-04-06, Excel serial 43927) of admissions data
# for all health organisations, replicating the pattern of the previous day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 43927

$names = @(
    "Akershus universitetssykehus HF",
    "Diakonhjemmet Sykehus",
    "Finnmarkssykehuset HF",
    "Haraldsplass Diakonale Sykehus",
    "Helgelandssykehuset HF",
    "Helse Bergen HF",
    "Helse Fonna HF",
    "Helse Førde HF",
    "Helse Møre og Romsdal",
    "Helse Nord-Trøndelag",
    "Helse Stavanger HF",
    "Lovisenberg Diakonale Sykehus",
    "Nordlandssykehuset HF",
    "Oslo universitetssykehus HF",
    "Sørlandet sykehus HF",
    "St. Olavs hospital",
    "Sunnaas sykehus HF",
    "Sykehuset Innlandet HF",
    "Sykehuset i Vestfold HF",
    "Sykehuset Østfold HF",
    "Sykehuset Telemark HF",
    "Universitetssykehuset Nord-Norge HF",
    "Vestre Viken HF"
)

$vals = @(55, 11, 2, 4, 0, 11, 6, 2, 8, 5, 10, 14, 6, 66, 12, 12, 0, 15, 8, 20, 6, 10, 33)

$startRow = 654

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newDate
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $vals[$i]
}
